$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 158 - this shifts rows 158..278 down to 159..279
$ws.Rows.Item(158).Insert()

# Copy the style/values that used to live in row 158 (now at row 159) into new row 158
# Columns that stay the same as the previous occupant of row 158 (now row 159)
$ws.Cells.Item(158, 1).Value = 3                                  # A - Mercado ID
$ws.Cells.Item(158, 2).Value = "Femacal de La Calera"             # B - Mercado
$ws.Cells.Item(158, 3).Value = "Coquimbo"                         # C - Region
$ws.Cells.Item(158, 4).Value = 44574                              # D - Fecha
$ws.Cells.Item(158, 5).Value = 5                                  # E - Codreg
$ws.Cells.Item(158, 6).Value = 100112040                          # F - Categoria ID
$ws.Cells.Item(158, 7).Value = "Cilantro"                         # G - Categoria
$ws.Cells.Item(158, 8).Value = "Sin especificar"                  # H - Variedad
$ws.Cells.Item(158, 9).Value = "Primera"                          # I - Calidad
$ws.Cells.Item(158, 10).Value = 125                                # J - Volumen
$ws.Cells.Item(158, 11).Value = 5000                               # K - Precio minimo
$ws.Cells.Item(158, 12).Value = 5000                               # L - Precio maximo
$ws.Cells.Item(158, 13).Value = 5000                               # M - Precio promedio ponderado
$ws.Cells.Item(158, 14).Value = "`$/docena de atados (3 kilos)"    # N - Unidad de comercializacion
$ws.Cells.Item(158, 15).Value = "Provincia de Quillota"            # O - Origen
$ws.Cells.Item(158, 16).Value = 1667                               # P - Precio $/Kg
$ws.Cells.Item(158, 17).Value = 3                                  # Q - Kg o Unidades
$ws.Cells.Item(158, 18).Value = "Hortaliza"                        # R - Clasificacion
